$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10-35 shift down to 11-36,
# and all formatting (including the date-formatted column D) shifts with them.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new data record.
$ws.Cells.Item(10, 1).Value2 = 6
$ws.Cells.Item(10, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(10, 3).Value2 = "Metropolitana"
$ws.Cells.Item(10, 4).Value2 = 44648
$ws.Cells.Item(10, 5).Value2 = 13
$ws.Cells.Item(10, 6).Value2 = "Fruta"
$ws.Cells.Item(10, 7).Value2 = 100102
$ws.Cells.Item(10, 8).Value2 = "Cítricos"
$ws.Cells.Item(10, 9).Value2 = 100102006
$ws.Cells.Item(10, 10).Value2 = "Pomelo"
$ws.Cells.Item(10, 11).Value2 = "Start Ruby"
$ws.Cells.Item(10, 12).Value2 = "Primera"
$ws.Cells.Item(10, 13).Value2 = 15
$ws.Cells.Item(10, 14).Value2 = 180000
$ws.Cells.Item(10, 15).Value2 = 180000
$ws.Cells.Item(10, 16).Value2 = 180000
$ws.Cells.Item(10, 17).Value2 = "$/bins (350 kilos)"
$ws.Cells.Item(10, 18).Value2 = "Región Metropolitana"
$ws.Cells.Item(10, 19).Value2 = 514
$ws.Cells.Item(10, 20).Value2 = 350
